$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.685.39"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").Value = "2.025.33"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'227.05"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").Value = "'0.604"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("D7").Value = "'59.72"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.374"
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("D10").Value = "'0.0824"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "2.327.51"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").Value = "'14.36"
$ws.Range("E13").Value = "  -3.00%  "
$ws.Range("D14").Value = "'20.96"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").Value = "'0.756"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "'5.15"
$ws.Range("E16").Value = "  -2.95%  "
$ws.Range("D17").Value = "2.036.92"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "37.614.59"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").Value = "'69.31"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").Value = "'5.89"
$ws.Range("E20").Value = "  -6.81%  "
$ws.Range("D21").Value = "0.0₃0821"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").Value = "'223.02"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "'2.37"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("D25").Value = "'2.24"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").Value = "'167.34"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "'9.21"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("E28").Value = "  -3.43%  "
$ws.Range("D29").Value = "'18.74"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("E30").Value = "  -5.76%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  +8.17%  "
$ws.Range("D33").Value = "'4.37"
$ws.Range("E33").Value = "  -4.11%  "
$ws.Range("D34").Value = "'0.0601"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").Value = "'4.45"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").Value = "'6.31"
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("D37").Value = "'2.28"
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "'17.77"
$ws.Range("E40").Value = "  +3.86%  "
$ws.Range("D41").Value = "1.535.73"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").Value = "'95.22"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("D45").Value = "'0.0909"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "'2.96"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("D50").Value = "'7.08"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").Value = "2.218.42"
$ws.Range("E51").Value = "  -1.48%  "
